# The match/odds data that used to be on row 11 (id 9) actually belongs to
# the match on row 12 (id 10), and vice versa. Swap the full data rows
# (everything except the fixed id/Div/Date/FTR columns A, C, D, K, which are
# identical for both matches) between rows 11 and 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row11 = $ws.Range("B11:AD11").Value()
$row12 = $ws.Range("B12:AD12").Value()

$ws.Range("B11:AD11").Value = $row12
$ws.Range("B12:AD12").Value = $row11
